$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates. Column D (Price) values are forced to Text format
# ("@") before assignment so that numeric-looking strings (e.g. "485.22",
# "61.60") are preserved verbatim as text instead of being auto-converted
# to numbers by Excel (which would also strip significant trailing zeros).

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "68.317.77"
$ws.Cells.Item(2, 5).Value = "  -0.56%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.903.55"
$ws.Cells.Item(3, 5).Value = "  -0.44%  "
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "485.22"
$ws.Cells.Item(5, 5).Value = "  +0.95%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "145.87"
$ws.Cells.Item(6, 5).Value = "  +0.82%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.622"
$ws.Cells.Item(7, 5).Value = "  +0.30%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.998"
$ws.Cells.Item(8, 5).Value = "  +0.04%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.741"
$ws.Cells.Item(9, 5).Value = "  +2.67%  "
$ws.Cells.Item(10, 5).Value = "  +7.77%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0000352"
$ws.Cells.Item(11, 5).Value = "  +0.19%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "42.94"
$ws.Cells.Item(12, 5).Value = "  +1.01%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.46"
$ws.Cells.Item(13, 5).Value = "  -0.16%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.516.31"
$ws.Cells.Item(14, 5).Value = "  -0.69%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.919.89"
$ws.Cells.Item(15, 5).Value = "  -0.87%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.22"
$ws.Cells.Item(16, 5).Value = "  -2.40%  "
$ws.Cells.Item(17, 5).Value = "  -0.56%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "19.99"
$ws.Cells.Item(18, 5).Value = "  +1.80%  "
$ws.Cells.Item(19, 5).Value = "  +0.96%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "68.422.42"
$ws.Cells.Item(20, 5).Value = "  -0.43%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "431.07"
$ws.Cells.Item(21, 5).Value = "  -0.29%  "
$ws.Cells.Item(22, 5).Value = "  +7.13%  "
$ws.Cells.Item(23, 5).Value = "  +1.10%  "
$ws.Cells.Item(24, 2).Value = "RenderToken"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "12.35"
$ws.Cells.Item(24, 5).Value = "  +21.56%  "
$ws.Cells.Item(25, 2).Value = "Litecoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "89.32"
$ws.Cells.Item(25, 5).Value = "  +2.32%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.71"
$ws.Cells.Item(26, 5).Value = "  +3.94%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.95"
$ws.Cells.Item(27, 5).Value = "  -5.98%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "37.31"
$ws.Cells.Item(28, 5).Value = "  -2.04%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.69"
$ws.Cells.Item(29, 5).Value = "  -2.08%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "717.28"
$ws.Cells.Item(30, 5).Value = "  +1.86%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.39"
$ws.Cells.Item(31, 5).Value = "  +1.10%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.131"
$ws.Cells.Item(32, 5).Value = "  +0.49%  "
$ws.Cells.Item(33, 5).Value = "  +2.09%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0₃0896"
$ws.Cells.Item(34, 5).Value = "  -1.29%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "61.60"
$ws.Cells.Item(35, 5).Value = "  +4.98%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "6.07"
$ws.Cells.Item(36, 5).Value = "  +8.50%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "40.75"
$ws.Cells.Item(37, 5).Value = "  -1.34%  "
$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.404"
$ws.Cells.Item(38, 5).Value = "  +19.16%  "
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.147"
$ws.Cells.Item(39, 5).Value = "  -2.73%  "
$ws.Cells.Item(40, 5).Value = "  +0.13%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0495"
$ws.Cells.Item(41, 5).Value = "  +4.91%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.98"
$ws.Cells.Item(42, 5).Value = "  +9.37%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.09"
$ws.Cells.Item(43, 5).Value = "  +3.74%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "3.03"
$ws.Cells.Item(44, 5).Value = "  -0.72%  "
$ws.Cells.Item(45, 5).Value = "  +0.56%  "
$ws.Cells.Item(46, 5).Value = "  +0.06%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.33"
$ws.Cells.Item(47, 5).Value = "  +5.94%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0₆0360"
$ws.Cells.Item(48, 5).Value = "  +24.80%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.37"
$ws.Cells.Item(49, 5).Value = "  -1.16%  "
$ws.Cells.Item(50, 5).Value = "  -2.39%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "143.75"
$ws.Cells.Item(51, 5).Value = "  -2.28%  "
